# Update generated output (想去人数 / F-column) counters and add one new
# show ("广州·三重唱Ohashi Trio（大桥トリ才） 2024年巡演") to the 演出 sheet
# (and its duplicate row inside the combined 全部类型 sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. 展览 ("exhibitions") sheet - plain counter bumps, no structural change
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    7  = 759
    8  = 261
    10 = 49
    12 = 221
    13 = 90
    14 = 901
    16 = 2023
    17 = 514
    18 = 8134
    19 = 666
    20 = 525
    21 = 69
    23 = 23
    24 = 235
    25 = 141
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value2 = $expoUpdates[$row]
}

# ---------------------------------------------------------------------
# 2. 演出 ("shows") sheet - insert a new row for the Ohashi Trio show
#    right above the existing "变形金刚..." row (old row 19).
# ---------------------------------------------------------------------
$wsShows = $wb.Worksheets.Item("演出")

$wsShows.Rows.Item(19).Copy()
$wsShows.Rows.Item(19).Insert()

# Re-apply the exact cell formatting from the row that is now beneath the
# freshly inserted one (A:I only, so we don't bleed formats into unused
# columns / rows).
$wsShows.Range("A20:I20").Copy()
$wsShows.Range("A19:I19").PasteSpecial(-4122)

$wsShows.Cells.Item(19, 1).Value2 = 18
$wsShows.Cells.Item(19, 2).NumberFormat = "@"
$wsShows.Cells.Item(19, 2).Value2 = "2024-11-27"
$wsShows.Cells.Item(19, 2).Style = "Normal"
$wsShows.Cells.Item(19, 3).Value2 = "广州·三重唱Ohashi Trio（大桥トリ才） 2024年巡演"
$wsShows.Cells.Item(19, 4).Value2 = "人民北路875号（广州市少年宫内） 广州蓓蕾剧院"
$wsShows.Cells.Item(19, 5).Value2 = "2024.11.27 19:30-11.27 21:00"
$wsShows.Cells.Item(19, 6).Value2 = 0
$wsShows.Cells.Item(19, 7).Value2 = 380
$wsShows.Cells.Item(19, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=91847"
$wsShows.Cells.Item(19, 9).Value2 = "//i2.hdslb.com/bfs/openplatform/202409/ggAAQH8D1725369168304.jpeg"

# The shifted rows keep their original content, but the running index in
# column A needs to stay in lock-step with (row number - 1).
$wsShows.Cells.Item(20, 1).Value2 = 19
$wsShows.Cells.Item(21, 1).Value2 = 20
$wsShows.Cells.Item(22, 1).Value2 = 21

# ---------------------------------------------------------------------
# 3. 本地生活 ("local life") sheet - plain counter bumps
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Cells.Item(2, 6).Value2 = 5571
$wsLocal.Cells.Item(4, 6).Value2 = 399

# ---------------------------------------------------------------------
# 4. 全部类型 ("all types") sheet - combination of the 展览 counter bumps
#    (on the rows that mirror 展览) plus the same new Ohashi Trio row
#    insertion that 演出 received (mirrored at row 46 there).
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Cells.Item(3, 6).Value2 = 5571
$wsAll.Cells.Item(5, 6).Value2 = 399

$allUpdates = @{
    13 = 759
    15 = 261
    18 = 49
    20 = 221
    22 = 90
    24 = 901
    28 = 2023
    29 = 514
    30 = 8134
    33 = 666
    34 = 525
    35 = 69
    38 = 23
    40 = 235
    42 = 141
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value2 = $allUpdates[$row]
}

$wsAll.Rows.Item(46).Copy()
$wsAll.Rows.Item(46).Insert()

$wsAll.Range("A47:I47").Copy()
$wsAll.Range("A46:I46").PasteSpecial(-4122)

$wsAll.Cells.Item(46, 1).Value2 = 45
$wsAll.Cells.Item(46, 2).NumberFormat = "@"
$wsAll.Cells.Item(46, 2).Value2 = "2024-11-27"
$wsAll.Cells.Item(46, 2).Style = "Normal"
$wsAll.Cells.Item(46, 3).Value2 = "广州·三重唱Ohashi Trio（大桥トリ才） 2024年巡演"
$wsAll.Cells.Item(46, 4).Value2 = "人民北路875号（广州市少年宫内） 广州蓓蕾剧院"
$wsAll.Cells.Item(46, 5).Value2 = "2024.11.27 19:30-11.27 21:00"
$wsAll.Cells.Item(46, 6).Value2 = 0
$wsAll.Cells.Item(46, 7).Value2 = 380
$wsAll.Cells.Item(46, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=91847"
$wsAll.Cells.Item(46, 9).Value2 = "//i2.hdslb.com/bfs/openplatform/202409/ggAAQH8D1725369168304.jpeg"

$wsAll.Cells.Item(47, 1).Value2 = 46
$wsAll.Cells.Item(48, 1).Value2 = 47
$wsAll.Cells.Item(49, 1).Value2 = 48

Write-Output "Edits applied."
